# Rename "Simulation score" to "Performance score"
$wb = $excel.ActiveWorkbook

# 1. Rename the worksheet tab. Excel automatically updates in-cell formula
#    references (e.g. 'Simulation score'!AC22 -> 'Performance score'!AC22)
#    when a sheet is renamed.
$ws1 = $wb.Worksheets.Item("Simulation score")
$ws1.Name = "Performance score"

# 2. Update the label cell on the (now renamed) score sheet that spelled out
#    the total score heading.
$ws1.Range("AB20").Value = "Total performancen scores"

# 3. Update the "Total score" sheet labels that referenced the old sheet name.
$ws3 = $wb.Worksheets.Item("Total score")
$ws3.Range("B2").Value = "Performance score"
$ws3.Range("B5").Value = "Performance score"

# Column B on "Total score" auto-sizes to fit the new (longer) label text.
$ws3.Columns.Item(2).ColumnWidth = 17

# 4. The chart on the renamed sheet keeps its own series formulas, which are
#    not automatically retargeted when the sheet is renamed, so update them
#    explicitly.
$co = $ws1.ChartObjects().Item(1)
$chart = $co.Chart
$cols = @("Q", "R", "S")
for ($i = 1; $i -le $chart.SeriesCollection().Count; $i++) {
    $s = $chart.SeriesCollection().Item($i)
    $col = $cols[$i - 1]
    $s.Formula = "=SERIES('Performance score'!`$" + $col + "`$1,,'Performance score'!`$" + $col + "`$2:`$" + $col + "`$49," + $i + ")"
}
